$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $row = $cell.Row
    $styleSource = $ws.Cells.Item($row, 3)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $styleSource.Style
}

$ws.Range("D2").Value = "37.338.89"
$ws.Range("E2").Value = "  +1.01%  "

$ws.Range("D3").Value = "2.031.12"
$ws.Range("E3").Value = "  +0.59%  "

$ws.Range("E4").Value = "  +0.23%  "

Set-TextValue $ws.Range("D5") "229.32"
$ws.Range("E5").Value = "  +1.19%  "

$ws.Range("E6").Value = "  +0.75%  "

$ws.Range("E7").Value = "  +0.02%  "

Set-TextValue $ws.Range("D8") "56.22"
$ws.Range("E8").Value = "  +2.23%  "

$ws.Range("E9").Value = "  -0.38%  "

Set-TextValue $ws.Range("D10") "0.0784"
$ws.Range("E10").Value = "  -1.05%  "

$ws.Range("E11").Value = "  -2.07%  "

$ws.Range("D12").Value = "2.335.14"
$ws.Range("E12").Value = "  +0.70%  "

Set-TextValue $ws.Range("D13") "14.34"
$ws.Range("E13").Value = "  +0.23%  "

Set-TextValue $ws.Range("D14") "20.31"
$ws.Range("E14").Value = "  -1.32%  "

Set-TextValue $ws.Range("D15") "0.742"
$ws.Range("E15").Value = "  -0.26%  "

$ws.Range("E16").Value = "  +0.94%  "

$ws.Range("D17").Value = "2.030.43"
$ws.Range("E17").Value = "  +0.42%  "

$ws.Range("D18").Value = "37.265.26"
$ws.Range("E18").Value = "  +1.00%  "

Set-TextValue $ws.Range("D19") "6.19"
$ws.Range("E19").Value = "  +1.65%  "

Set-TextValue $ws.Range("D20") "68.96"
$ws.Range("E20").Value = "  +0.22%  "

$ws.Range("D21").Value = "0.0₃0819"
$ws.Range("E21").Value = "  -1.01%  "

Set-TextValue $ws.Range("D22") "223.60"
$ws.Range("E22").Value = "  -1.36%  "

$ws.Range("E23").Value = "  +0.03%  "

$ws.Range("E24").Value = "  +1.95%  "

$ws.Range("E25").Value = "  -1.05%  "

Set-TextValue $ws.Range("D26") "164.53"
$ws.Range("E26").Value = "  -1.47%  "

Set-TextValue $ws.Range("D27") "9.10"
$ws.Range("E27").Value = "  -1.91%  "

Set-TextValue $ws.Range("D28") "0.133"
$ws.Range("E28").Value = "  +4.55%  "

Set-TextValue $ws.Range("D29") "18.74"
$ws.Range("E29").Value = "  -0.12%  "

Set-TextValue $ws.Range("D30") "1.32"
$ws.Range("E30").Value = "  -1.69%  "

$ws.Range("E31").Value = "  +0.46%  "

$ws.Range("E32").Value = "  -0.40%  "

Set-TextValue $ws.Range("D33") "0.0607"
$ws.Range("E33").Value = "  -0.86%  "

$ws.Range("E34").Value = "  +9.60%  "

Set-TextValue $ws.Range("D35") "4.45"
$ws.Range("E35").Value = "  -0.14%  "

$ws.Range("E36").Value = "  -1.99%  "

Set-TextValue $ws.Range("D37") "3.24"
$ws.Range("E37").Value = "  +1.80%  "

$ws.Range("E38").Value = "  +0.18%  "

Set-TextValue $ws.Range("D39") "5.62"
$ws.Range("E39").Value = "  +4.19%  "

$ws.Range("D40").Value = "1.471.22"
$ws.Range("E40").Value = "  -1.32%  "

$ws.Range("E41").Value = "  -2.43%  "

Set-TextValue $ws.Range("D42") "4.33"
$ws.Range("E42").Value = "  +16.13%  "

$ws.Range("E43").Value = "  +0.35%  "

$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D44") "0.0921"
$ws.Range("E44").Value = "  -0.91%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D45") "94.29"
$ws.Range("E45").Value = "  -0.81%  "

Set-TextValue $ws.Range("D46") "16.26"
$ws.Range("E46").Value = "  -4.76%  "

$ws.Range("E47").Value = "  -2.37%  "

$ws.Range("E48").Value = "  +0.70%  "

Set-TextValue $ws.Range("D49") "7.11"
$ws.Range("E49").Value = "  -2.38%  "

$ws.Range("E50").Value = "  +0.97%  "

$ws.Range("D51").Value = "2.219.25"
$ws.Range("E51").Value = "  +0.48%  "
